# Actualización automática 2025-06-27 17:25:45
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("E6").Value = 199.74
$wsVentasGrupo.Range("M6").Value = 73.34
$wsVentasGrupo.Range("M29").Value = 535.85
$wsVentasGrupo.Range("M55").Value = "10 de 53"

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F6").Value = 1516.28
$wsVentaMensual.Range("F29").Value = 1079.23
$wsVentaMensual.Range("F55").Value = 88195.77

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D2").Value = 5502.24
$wsCumplimiento.Range("E2").Value = 317.7600000000002
$wsCumplimiento.Range("F2").Value = 0.9454020618556701

$wsCumplimiento.Range("D4").Value = 871.86
$wsCumplimiento.Range("E4").Value = -225.86
$wsCumplimiento.Range("F4").Value = 1.349628482972136

$wsCumplimiento.Range("D15").Value = 11324.01
$wsCumplimiento.Range("E15").Value = 4365.99
$wsCumplimiento.Range("F15").Value = 0.721734225621415

$wsCumplimiento.Range("D16").Value = 43477.13
$wsCumplimiento.Range("E16").Value = 2268.559000000001
$wsCumplimiento.Range("F16").Value = 0.9504093380252727

$wsCumplimiento.Range("D19").Value = 99631.23
$wsCumplimiento.Range("E19").Value = -8667.900999999998
$wsCumplimiento.Range("F19").Value = 1.095290059140206
